# Humedad y Anexo listo
#
# 1) Rebuild the "{#tractoras}...{/tractoras}" block: split the old
#    single "Tractora con matricula {{matriculas}} y con remolque
#    {{remolque}}" paragraph (which illegally mixed the loop-open tag
#    into its own text) into three clean paragraphs, switch the
#    placeholders from {{double}} to {single} brace syntax and rename
#    them to matriculaTractora / matriculaRemolque, and add two blank
#    paragraphs after the loop closes.
# 2) Tidy the "{{fecha}}" placeholder down to "{fecha}".
# 3) Turn the second blank paragraph after the fecha line into an
#    (empty) underlined paragraph, mirroring the one above the
#    "PREVISTOS..." paragraph.

$d = $word.ActiveDocument

function Get-ParagraphIndexContaining($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# --- 1) tractoras loop block -------------------------------------------------

$loopStartIdx = Get-ParagraphIndexContaining $d "tractoras"
$loopEndIdx = $loopStartIdx + 1

$loopRange = $d.Range($d.Paragraphs.Item($loopStartIdx).Range.Start, `
                       $d.Paragraphs.Item($loopEndIdx).Range.End)

$tractorasXml = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p><w:r><w:t>{#tractoras}</w:t></w:r></w:p><w:p><w:r><w:t>Tractora con matr&#237;cula {</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>matriculaTractora</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">} y con remolque </w:t></w:r><w:r><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>matriculaRemolque</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r></w:p><w:p><w:r><w:t>{/tractoras}</w:t></w:r></w:p><w:p/><w:p/></pkg:xmlData>
'@

$loopRange.InsertXML($tractorasXml)

# --- 2) {{fecha}} -> {fecha} --------------------------------------------------

$fechaIdx = Get-ParagraphIndexContaining $d "fecha"
$fechaRange = $d.Range($d.Paragraphs.Item($fechaIdx).Range.Start, `
                        $d.Paragraphs.Item($fechaIdx).Range.End)

$fechaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">En Fuentes de </w:t></w:r><w:r><w:t>Andaluc&#237;a</w:t></w:r><w:r><w:t xml:space="preserve"> a </w:t></w:r><w:r><w:t>{</w:t></w:r><w:r><w:t>fecha</w:t></w:r><w:r><w:t>}</w:t></w:r></w:p>'

$fechaRange.InsertXML($fechaXml)

# --- 3) second blank paragraph after the fecha line becomes underlined ------

$underlineIdx = $fechaIdx + 2
$underlineRange = $d.Range($d.Paragraphs.Item($underlineIdx).Range.Start, `
                            $d.Paragraphs.Item($underlineIdx).Range.End)

$underlineXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr></w:p>'

$underlineRange.InsertXML($underlineXml)
